# Update "想去人数" (interested-count) values on the "展览" and "全部类型"
# sheets to reflect the latest scrape at commit 456a3b4:
#   F2: 203 -> 202
#   F3: 154 -> 155

$wb = $excel.ActiveWorkbook

foreach ($sheetName in @("展览", "全部类型")) {
    $ws = $wb.Worksheets.Item($sheetName)
    $ws.Range("F2").Value = 202
    $ws.Range("F3").Value = 155
}
